$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Anonymize the "fedcore" header label.
$ws1.Range("C2").Value = "approach"

# Normalize "-0" artifacts to plain "0".
$ws1.Range("D4").Value = 0
$ws1.Range("D10").Value = 0
$ws1.Range("D12").Value = 0

# Give the merged title row (B1:D1) a real border outline: the two inner
# cells of the merge (C1 and D1) pick up their own border-only styles
# (top+bottom for C1, top+bottom+right for D1), while B1 keeps its
# existing full-box style untouched.
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$d1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$d1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous

$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$c1.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous

# ---------------------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Anonymize both "fedcore" header labels (B1:D1 and E1:G1 merged blocks).
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Same border treatment as sheet 1, for both merged title ranges.
$d1b = $ws2.Range("D1")
$d1b.Style = "Normal"
$d1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$d1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$d1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous

$c1b = $ws2.Range("C1")
$c1b.Style = "Normal"
$c1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$c1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous

$g1b = $ws2.Range("G1")
$g1b.Style = "Normal"
$g1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$g1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$g1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous

$f1b = $ws2.Range("F1")
$f1b.Style = "Normal"
$f1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$f1b.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous

# Drop the stray empty placeholder cell at G5.
$ws2.Range("G5").ClearContents()
